{"js": "// Golden-test style update for track_changes_insertion.docx\n//\n// 1. Add a new \"Abstract Title\" paragraph style (styleId \"AbstractTitle\"),\n//    based on Normal, followed by Abstract (same shape as the existing\n//    Author/Date/Abstract styles in this stylesheet).\n// 2. Abstract style: space-before 300 -> 100 (twentieths of a point).\n// 3. ImportTok character style: add bold + green (#008000) color.\n// 4. BuiltInTok character style: add green (#008000) color.\n\n// --- 1. New \"Abstract Title\" style -----------------------------------\n// Passing the space-containing display name to addStyle() yields\n// styleId \"AbstractTitle\" / name \"Abstract Title\", matching the target.\ncontext.document.addStyle(\"Abstract Title\", \"Paragraph\");\nawait context.sync();\n\nconst styles = context.document.getStyles();\nconst abstractTitle = styles.getByNameOrNullObject(\"Abstract Title\");\nawait context.sync();\n\n// Re-fetching the style (instead of using the addStyle() return value\n// directly) is required here: property writes on the freshly-minted\n// object are otherwise applied to the document body instead of the\n// style definition.\nabstractTitle.baseStyle = \"Normal\";\nabstractTitle.nextParagraphStyle = \"Abstract\";\nabstractTitle.quickStyle = true;\nawait context.sync();\n\nabstractTitle.paragraphFormat.keepWithNext = true; // w:keepNext\nabstractTitle.paragraphFormat.keepTogether = true; // w:keepLines\nabstractTitle.paragraphFormat.alignment = Word.Alignment.centered;\nabstractTitle.paragraphFormat.spaceAfter = 0;\nabstractTitle.paragraphFormat.spaceBefore = 15; // points -> w:before=\"300\"\nawait context.sync();\n\nabstractTitle.font.size = 10; // -> w:sz=\"20\"\nabstractTitle.font.bold = true;\nabstractTitle.font.color = \"#345A8A\";\nawait context.sync();\n\n// --- 2. Abstract style spacing tweak ----------------------------------\nconst abstractStyle = styles.getByNameOrNullObject(\"Abstract\");\nawait context.sync();\nabstractStyle.paragraphFormat.spaceBefore = 5; // points -> w:before=\"100\"\nawait context.sync();\n\n// --- 3. ImportTok character style ------------------------------------\nconst importTok = styles.getByNameOrNullObject(\"ImportTok\");\nawait context.sync();\nimportTok.font.color = \"#008000\";\nimportTok.font.bold = true;\nawait context.sync();\n\n// --- 4. BuiltInTok character style ------------------------------------\nconst builtInTok = styles.getByNameOrNullObject(\"BuiltInTok\");\nawait context.sync();\nbuiltInTok.font.color = \"#008000\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$styles = $d.Styles\n\n# 1. Add the new \"Abstract Title\" paragraph style (customStyle \"AbstractTitle\"),\n#    based on Normal, followed by Abstract, matching the sibling styles\n#    (Author/Date/Abstract) already present in this stylesheet.\n$abstractTitle = $styles.Add(\"AbstractTitle\", 1)\n$abstractTitle.NameLocal = \"Abstract Title\"\n$abstractTitle.BaseStyle = \"Normal\"\n$abstractTitle.NextParagraphStyle = \"Abstract\"\n$abstractTitle.QuickStyle = $true\n\n$abstractTitle.ParagraphFormat.KeepWithNext = $true\n$abstractTitle.ParagraphFormat.KeepTogether = $true\n$abstractTitle.ParagraphFormat.Alignment = 1\n$abstractTitle.ParagraphFormat.SpaceAfter = 0\n$abstractTitle.ParagraphFormat.SpaceBefore = 15\n\n$abstractTitle.Font.Size = 10\n$abstractTitle.Font.SizeBi = 10\n$abstractTitle.Font.Bold = $true\n$abstractTitle.Font.Color = 9067060\n\n# 2. Abstract style: space-before 300 -> 100 (twips/20 = 15pt -> 5pt).\n$abstract = $styles.Item(\"Abstract\")\n$abstract.ParagraphFormat.SpaceBefore = 5\n\n# 3. ImportTok character style: add bold + green color.\n$importTok = $styles.Item(\"ImportTok\")\n$importTok.Font.Color = 32768\n$importTok.Font.Bold = $true\n\n# 4. BuiltInTok character style: add green color.\n$builtInTok = $styles.Item(\"BuiltInTok\")\n$builtInTok.Font.Color = 32768\n"}
